# "view and add student in college correctedd"
#
# The "event_id" column (F1 header on Sheet1) is removed/cleared - the
# student roster no longer tracks an event id for each student. Clearing
# the cell's contents also drops the now-unused "event_id" entry from the
# shared-strings table on save, and re-selects F1:F3 (the column that was
# just edited) the way Excel leaves the just-edited range selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "event_id" header text from F1 (keep its formatting/style).
$ws.Range("F1").ClearContents()

# Leave the edited column selected, as reflected in the saved sheet view.
$ws.Range("F1:F3").Select()
